$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1. "Metas e Viabilidade" + " - Byte INC." used to be two runs with identical
#    formatting; they become a single run with the same combined text.
$d.Content.Find.Execute("Metas e Viabilidade – Byte INC.", $true, $false, $false, $false, $false, $true, 1, $false, "Metas e Viabilidade – Byte INC.", 2) | Out-Null

# 2. Likewise, the "purpose" sentence's two runs merge into one run.
$d.Content.Find.Execute("Este documento tem como propósito descrever metas e apresentar um estudo de viabilidade para o projeto EveRemind.", $true, $false, $false, $false, $false, $true, 1, $false, "Este documento tem como propósito descrever metas e apresentar um estudo de viabilidade para o projeto EveRemind.", 2) | Out-Null

# 3. "Entregar um produto com qualidade conforme os requisitos elicitados;"
#    splits into three runs and "elicitados" becomes "licitados".
$full = $d.Content.Text
$oldText3 = "Entregar um produto com qualidade conforme os requisitos elicitados;"
$start3 = $full.IndexOf($oldText3)
if ($start3 -lt 0) { throw "target text for change #3 not found" }
$rng3 = $d.Range($start3, $start3 + $oldText3.Length)
$xml3 = "<w:p $wNs>" +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Entregar um produto com qualidade </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">conforme os requisitos </w:t></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>licitados;</w:t></w:r>' +
  '</w:p>'
$rng3.InsertXML($xml3) | Out-Null

# 4. "A viabilidade ... com uma único motivo de cumprir e atingit as metas do
#    projeto:" becomes "... com um único motivo de cumprir e atingir as metas
#    do projeto:", with the _GoBack bookmark relocated right after "atingir"
#    (i.e. where the edit actually happened) instead of right before "motivo".
$full = $d.Content.Text
$oldText4 = "A viabilidade de um projeto pode ser pensada em três áreas com uma único motivo de cumprir e atingit as metas do projeto:"
$start4 = $full.IndexOf($oldText4)
if ($start4 -lt 0) { throw "target text for change #4 not found" }
$rng4 = $d.Range($start4, $start4 + $oldText4.Length)
$xml4 = "<w:p $wNs>" +
  '<w:r><w:t>A viabilidade de um projeto pode ser pensada em três áreas</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> com um</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> ú</w:t></w:r>' +
  '<w:r><w:t>nico motivo de cumprir e atingir</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t xml:space="preserve"> as metas do projeto</w:t></w:r>' +
  '<w:r><w:t>:</w:t></w:r>' +
  '</w:p>'
$rng4.InsertXML($xml4) | Out-Null
